$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$bfValues = New-Object 'object[,]' 24,5
$bfValues[0,0] = 2.826683192429698
$bfValues[0,1] = 0.364818345506535
$bfValues[0,2] = 0.009022407014562361
$bfValues[0,3] = 0.05043111356051849
$bfValues[0,4] = 4.157605467521165
$bfValues[1,0] = 2.736422786721562
$bfValues[1,1] = 0.3374530303535721
$bfValues[1,2] = 0.008450282533921438
$bfValues[1,3] = 0.05050730774579704
$bfValues[1,4] = 4.128671466142848
$bfValues[2,0] = 2.682727031003708
$bfValues[2,1] = 0.3208750531878479
$bfValues[2,2] = 0.008096173153653297
$bfValues[2,3] = 0.05055825357178495
$bfValues[2,4] = 4.112946417056065
$bfValues[3,0] = 2.661278921277699
$bfValues[3,1] = 0.3141752933429416
$bfValues[3,2] = 0.007951108720849476
$bfValues[3,3] = 0.05058006263721271
$bfValues[3,4] = 4.107049310231005
$bfValues[4,0] = 2.657743644695643
$bfValues[4,1] = 0.3130661641408778
$bfValues[4,2] = 0.007926973337678334
$bfValues[4,3] = 0.05058374737421756
$bfValues[4,4] = 4.106100899938511
$bfValues[5,0] = 2.682436020063619
$bfValues[5,1] = 0.3207844722416269
$bfValues[5,2] = 0.008094219917428092
$bfValues[5,3] = 0.05055854344986926
$bfValues[5,4] = 4.11286482053535
$bfValues[6,0] = 2.795203256877016
$bfValues[6,1] = 0.3553358682492274
$bfValues[6,2] = 0.008825692755099368
$bfValues[6,3] = 0.05045652236702669
$bfValues[6,4] = 4.147204202874363
$bfValues[7,0] = 3.030053990057752
$bfValues[7,1] = 0.4249002359604219
$bfValues[7,2] = 0.01024010721108581
$bfValues[7,3] = 0.05028941902577372
$bfValues[7,4] = 4.23083882366646
$bfValues[8,0] = 3.211029415870655
$bfValues[8,1] = 0.4771574687175644
$bfValues[8,2] = 0.01127045713545627
$bfValues[8,3] = 0.0501866576612664
$bfValues[8,4] = 4.302375552408307
$bfValues[9,0] = 3.295208059933657
$bfValues[9,1] = 0.5011903157224538
$bfValues[9,2] = 0.01173800890620313
$bfValues[9,3] = 0.05014423755449759
$bfValues[9,4] = 4.337145307364437
$bfValues[10,0] = 3.327351673147632
$bfValues[10,1] = 0.5103291676619506
$bfValues[10,2] = 0.01191495070795412
$bfValues[10,3] = 0.05012879513059898
$bfValues[10,4] = 4.350634647706045
$bfValues[11,0] = 3.320417085226893
$bfValues[11,1] = 0.5083592460607633
$bfValues[11,2] = 0.01187684710344428
$bfValues[11,3] = 0.05013209332226654
$bfValues[11,4] = 4.34771508063082
$bfValues[12,0] = 3.297847179447785
$bfValues[12,1] = 0.5019414066056243
$bfValues[12,2] = 0.01175256788974721
$bfValues[12,3] = 0.05014295465367463
$bfValues[12,4] = 4.338248599914493
$bfValues[13,0] = 3.284057251091383
$bfValues[13,1] = 0.4980152799806774
$bfValues[13,2] = 0.01167643066466439
$bfValues[13,3] = 0.05014968839634631
$bfValues[13,4] = 4.332492218494906
$bfValues[14,0] = 3.205565466644998
$bfValues[14,1] = 0.4755921670246153
$bfValues[14,2] = 0.01123988304239987
$bfValues[14,3] = 0.05018951688893714
$bfValues[14,4] = 4.300148299900457
$bfValues[15,0] = 3.157888101739218
$bfValues[15,1] = 0.4619035507003559
$bfValues[15,2] = 0.01097181851917739
$bfValues[15,3] = 0.05021505777606672
$bfValues[15,4] = 4.280878597272249
$bfValues[16,0] = 3.13063961897916
$bfValues[16,1] = 0.4540547265776809
$bfValues[16,2] = 0.01081752232342126
$bfValues[16,3] = 0.05023015550990084
$bfValues[16,4] = 4.270004670293503
$bfValues[17,0] = 3.121443663654816
$bfValues[17,1] = 0.4514014393768662
$bfValues[17,2] = 0.01076525922414362
$bfValues[17,3] = 0.05023533732753666
$bfValues[17,4] = 4.266358850947427
$bfValues[18,0] = 3.162945397016983
$bfValues[18,1] = 0.4633581864458165
$bfValues[18,2] = 0.01100036568516316
$bfValues[18,3] = 0.05021229675947581
$bfValues[18,4] = 4.282908189936137
$bfValues[19,0] = 3.304469256257562
$bfValues[19,1] = 0.503825441525521
$bfValues[19,2] = 0.01178907423478748
$bfValues[19,3] = 0.05013974756726314
$bfValues[19,4] = 4.341020354293022
$bfValues[20,0] = 3.39851984345006
$bfValues[20,1] = 0.5304957142211038
$bfValues[20,2] = 0.01230392696546545
$bfValues[20,3] = 0.05009595260102973
$bfValues[20,4] = 4.380882655891128
$bfValues[21,0] = 3.348180611550276
$bfValues[21,1] = 0.5162407074033126
$bfValues[21,2] = 0.01202917764272726
$bfValues[21,3] = 0.05011899586906865
$bfValues[21,4] = 4.359434298955449
$bfValues[22,0] = 3.160658489684067
$bfValues[22,1] = 0.4627004804382864
$bfValues[22,2] = 0.01098746008372586
$bfValues[22,3] = 0.05021354372565962
$bfValues[22,4] = 4.281989974366581
$bfValues[23,0] = 2.965045411251083
$bfValues[23,1] = 0.4058829077467578
$bfValues[23,2] = 0.009859264737006157
$bfValues[23,3] = 0.0503311052509659
$bfValues[23,4] = 4.206452866420989

$ijValues = New-Object 'object[,]' 24,2
$ijValues[0,0] = 2.446716918797108
$ijValues[0,1] = 0.1247199570505728
$ijValues[1,0] = 2.436935429613143
$ijValues[1,1] = 0.1252668101824508
$ijValues[2,0] = 2.43204908349783
$ijValues[2,1] = 0.1256277414358475
$ijValues[3,0] = 2.430338114348771
$ijValues[3,1] = 0.1257811514937242
$ijValues[4,0] = 2.430070899374982
$ijValues[4,1] = 0.1258070073108257
$ijValues[5,0] = 2.432024875711065
$ijValues[5,1] = 0.1256297847564962
$ijValues[6,0] = 2.443111137013148
$ijValues[6,1] = 0.1249032902056229
$ijValues[7,0] = 2.473794460819718
$ijValues[7,1] = 0.1236782477280762
$ijValues[8,0] = 2.501878032762718
$ijValues[8,1] = 0.1228998502986514
$ijValues[9,0] = 2.515876609817653
$ijValues[9,1] = 0.1225721423331834
$ijValues[10,0] = 2.52135491836404
$ijValues[10,1] = 0.1224518426537085
$ijValues[11,0] = 2.520167155542396
$ijValues[11,1] = 0.1224775824788153
$ijValues[12,0] = 2.51632375079862
$ijValues[12,1] = 0.1225621691152696
$ijValues[13,0] = 2.513992694172828
$ijValues[13,1] = 0.1226144752979721
$ijValues[14,0] = 2.500987931072174
$ijValues[14,1] = 0.1229217979255939
$ijValues[15,0] = 2.493324255758495
$ijValues[15,1] = 0.1231170902945955
$ijValues[16,0] = 2.489031359859169
$ijValues[16,1] = 0.1232319010554921
$ijValues[17,0] = 2.487597574498295
$ijValues[17,1] = 0.1232712005487731
$ijValues[18,0] = 2.494128149624203
$ijValues[18,1] = 0.1230960440187729
$ijValues[19,0] = 2.5174478270158
$ijValues[19,1] = 0.1225372209289688
$ijValues[20,0] = 2.53372300812228
$ijValues[20,1] = 0.1221941257473809
$ijValues[21,0] = 2.524941503606939
$ijValues[21,1] = 0.1223752168637624
$ijValues[22,0] = 2.493764357182172
$ijValues[22,1] = 0.1231055511456489
$ijValues[23,0] = 2.464527273448098
$ijValues[23,1] = 0.1239882895614688

$lValues = New-Object 'object[,]' 24,1
$lValues[0,0] = 0.4654450911102685
$lValues[1,0] = 0.4601784076898952
$lValues[2,0] = 0.4571847648224008
$lValues[3,0] = 0.4560252048454174
$lValues[4,0] = 0.4558363071262761
$lValues[5,0] = 0.4571688821454956
$lValues[6,0] = 0.463579292991426
$lValues[7,0] = 0.4780571534617906
$lValues[8,0] = 0.4898613003190206
$lValues[9,0] = 0.4954859969538035
$lValues[10,0] = 0.4976526515261099
$lValues[11,0] = 0.497184390497253
$lValues[12,0] = 0.4956635130472762
$lValues[13,0] = 0.494736713079007
$lValues[14,0] = 0.4894988467752341
$lValues[15,0] = 0.4863509098898788
$lValues[16,0] = 0.4845642921305569
$lValues[17,0] = 0.4839634934601236
$lValues[18,0] = 0.4866835297256671
$lValues[19,0] = 0.4961092350969807
$lValues[20,0] = 0.502483442498729
$lValues[21,0] = 0.4990618139454597
$lValues[22,0] = 0.4865330801217311
$lValues[23,0] = 0.4739358810965228

$ws.Range("B2:F25").Value = $bfValues
$ws.Range("I2:J25").Value = $ijValues
$ws.Range("L2:L25").Value = $lValues
